# Applies the numeric-value changes from the commit diff across all
# affected sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). All target
# cells are plain numeric literals (no formulas in this workbook).
$wb = $excel.ActiveWorkbook

# ===================== Sheet ALC =====================
$ws = $wb.Worksheets.Item("ALC")
# row 28 (G28=27772)
$ws.Range("H28").Value = 8794.579
$ws.Range("I28").Value = 317.8125
$ws.Range("J28").Value = 54004
$ws.Range("K28").Value = 317.8125
$ws.Range("L28").Value = 54004
$ws.Range("M28").Value = 167.1875
$ws.Range("N28").Value = -54974
# row 51 (G51=5486)
$ws.Range("H51").Value = 2450
$ws.Range("I51").Value = 1525
$ws.Range("J51").Value = 3375
$ws.Range("K51").Value = 1525
$ws.Range("L51").Value = 3375
$ws.Range("M51").Value = -1041
$ws.Range("N51").Value = -4343
# row 132 (G132=44049)
$ws.Range("H132").Value = 2930.4102
$ws.Range("I132").Value = 2785.75
$ws.Range("K132").Value = 8357.25
$ws.Range("M132").Value = -5827.25

# ===================== Sheet ARM =====================
$ws = $wb.Worksheets.Item("ARM")
# row 32 (G32=44147)
$ws.Range("H32").Value = 3086.05
$ws.Range("I32").Value = 3086.05
$ws.Range("K32").Value = 3086.05
$ws.Range("M32").Value = -2799.05
# row 97 (G97=19941)
$ws.Range("H97").Value = 1220.8823
$ws.Range("I97").Value = 987.9167
$ws.Range("J97").Value = 1780
$ws.Range("K97").Value = 987.9167
$ws.Range("L97").Value = 1780
$ws.Range("M97").Value = -491.9167
$ws.Range("N97").Value = -2772
# row 106 (G106=18679)
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = ""
# row 122 (G122=36168)
$ws.Range("H122").Value = 60564.06
$ws.Range("I122").Value = 68372.60000000001
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 205117.8
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -202667.8
$ws.Range("N122").Value = -10900
# row 132 (G132=43997)
$ws.Range("H132").Value = 1028259.06
$ws.Range("I132").Value = 1899.6492
$ws.Range("J132").Value = 4278397
$ws.Range("K132").Value = 5698.9476
$ws.Range("L132").Value = 12835191
$ws.Range("M132").Value = -3168.9476
$ws.Range("N132").Value = -12840251

# ===================== Sheet BSM =====================
$ws = $wb.Worksheets.Item("BSM")
# row 86 (G86=12526)
$ws.Range("H86").Value = 1817.8334
$ws.Range("I86").Value = 1817.0769
$ws.Range("J86").Value = 1821.1111
$ws.Range("K86").Value = 1817.0769
$ws.Range("L86").Value = 1821.1111
$ws.Range("M86").Value = -694.0769
$ws.Range("N86").Value = -4067.1111
# row 89 (G89=12526)
$ws.Range("H89").Value = 1817.8334
$ws.Range("I89").Value = 1817.0769
$ws.Range("J89").Value = 1821.1111
$ws.Range("K89").Value = 9085.3845
$ws.Range("L89").Value = 9105.5555
$ws.Range("M89").Value = -3469.3845
$ws.Range("N89").Value = -20337.5555
# row 94 (G94=19939)
$ws.Range("H94").Value = 436.65
$ws.Range("I94").Value = 437.35715
$ws.Range("J94").Value = 435
$ws.Range("K94").Value = 437.35715
$ws.Range("L94").Value = 435
$ws.Range("M94").Value = 13.64285000000001
$ws.Range("N94").Value = -1337
# row 134 (G134=43998)
$ws.Range("H134").Value = 2828.182
$ws.Range("I134").Value = 2735.838
$ws.Range("J134").Value = 3316.2856
$ws.Range("K134").Value = 8207.514000000001
$ws.Range("L134").Value = 9948.856800000001
$ws.Range("M134").Value = -5672.514000000001
$ws.Range("N134").Value = -15018.8568

# ===================== Sheet CRP =====================
$ws = $wb.Worksheets.Item("CRP")
# row 28 (G28=18348)
$ws.Range("H28").Value = 47637.2
$ws.Range("J28").Value = 47637.2
$ws.Range("L28").Value = 47637.2
$ws.Range("N28").Value = -48127.2
# row 31 (G31=44023)
$ws.Range("H31").Value = 4622.58
$ws.Range("I31").Value = 1521.0189
$ws.Range("J31").Value = 8120.085
$ws.Range("K31").Value = 1521.0189
$ws.Range("L31").Value = 8120.085
$ws.Range("M31").Value = -1226.0189
$ws.Range("N31").Value = -8710.084999999999
# row 34 (G34=44023)
$ws.Range("H34").Value = 4622.58
$ws.Range("I34").Value = 1521.0189
$ws.Range("J34").Value = 8120.085
$ws.Range("K34").Value = 1521.0189
$ws.Range("L34").Value = 8120.085
$ws.Range("M34").Value = -1319.0189
$ws.Range("N34").Value = -8524.084999999999
# row 43 (G43=18504)
$ws.Range("H43").Value = 94828.5
$ws.Range("J43").Value = 94828.5
$ws.Range("L43").Value = 94828.5
$ws.Range("N43").Value = -95196.5
# row 101 (G101=18504)
$ws.Range("H101").Value = 94828.5
$ws.Range("J101").Value = 94828.5
$ws.Range("L101").Value = 94828.5
$ws.Range("N101").Value = -101318.5
# row 134 (G134=44020)
$ws.Range("H134").Value = 5957036
$ws.Range("I134").Value = 8338333
$ws.Range("J134").Value = 3793.9167
$ws.Range("K134").Value = 25014999
$ws.Range("L134").Value = 11381.7501
$ws.Range("M134").Value = -25012464
$ws.Range("N134").Value = -16451.7501

# ===================== Sheet CUL =====================
$ws = $wb.Worksheets.Item("CUL")
# row 5 (G5=43974)
$ws.Range("H5").Value = 1265.6342
$ws.Range("I5").Value = 989.5
$ws.Range("J5").Value = 1744.2667
$ws.Range("K5").Value = 2968.5
$ws.Range("L5").Value = 5232.800099999999
$ws.Range("M5").Value = -2856.5
$ws.Range("N5").Value = -5456.800099999999
# row 94 (G94=19811)
$ws.Range("H94").Value = 2322.8333
$ws.Range("I94").Value = 383.33334
$ws.Range("J94").Value = 4262.3335
$ws.Range("K94").Value = 1150.00002
$ws.Range("L94").Value = 12787.0005
$ws.Range("M94").Value = -474.0000199999999
$ws.Range("N94").Value = -14139.0005
# row 113 (G113=27843)
$ws.Range("H113").Value = 630.2766
$ws.Range("I113").Value = 609.1539
$ws.Range("K113").Value = 1827.4617
$ws.Range("M113").Value = 342.5382999999999
# row 125 (G125=36043)
$ws.Range("H125").Value = 1125.4
$ws.Range("I125").Value = 617
$ws.Range("J125").Value = 1888
$ws.Range("K125").Value = 1851
$ws.Range("L125").Value = 5664
$ws.Range("M125").Value = 3069
$ws.Range("N125").Value = -15504
# row 134 (G134=44074)
$ws.Range("H134").Value = 5787.2036
$ws.Range("I134").Value = 2170.3333
$ws.Range("J134").Value = 8838.9375
$ws.Range("K134").Value = 6510.999899999999
$ws.Range("L134").Value = 26516.8125
$ws.Range("M134").Value = -1440.999899999999
$ws.Range("N134").Value = -36656.8125
# row 135 (G135=43974)
$ws.Range("H135").Value = 1265.6342
$ws.Range("I135").Value = 989.5
$ws.Range("J135").Value = 1744.2667
$ws.Range("K135").Value = 8905.5
$ws.Range("L135").Value = 15698.4003
$ws.Range("M135").Value = -6370.5
$ws.Range("N135").Value = -20768.4003
# row 137 (G137=44088)
$ws.Range("H137").Value = 36026.91
$ws.Range("I137").Value = 6253.6665
$ws.Range("K137").Value = 18760.9995
$ws.Range("M137").Value = -13660.9995
# row 138 (G138=44105)
$ws.Range("H138").Value = 2376.45
$ws.Range("I138").Value = 1196.0555
$ws.Range("J138").Value = 13000
$ws.Range("K138").Value = 3588.1665
$ws.Range("L138").Value = 39000
$ws.Range("M138").Value = 1551.8335
$ws.Range("N138").Value = -49280

# ===================== Sheet GSM =====================
$ws = $wb.Worksheets.Item("GSM")
# row 122 (G122=36182)
$ws.Range("H122").Value = 2732.5386
$ws.Range("I122").Value = 3760
$ws.Range("J122").Value = 2090.375
$ws.Range("K122").Value = 11280
$ws.Range("L122").Value = 6271.125
$ws.Range("M122").Value = -8830
$ws.Range("N122").Value = -11171.125
# row 123 (G123=34150)
$ws.Range("H123").Value = 8058.0835
$ws.Range("J123").Value = 9744.111000000001
$ws.Range("L123").Value = 9744.111000000001
$ws.Range("N123").Value = -14644.111
# row 132 (G132=44008)
$ws.Range("H132").Value = 27783062
$ws.Range("I132").Value = 41672736
$ws.Range("J132").Value = 3718.75
$ws.Range("K132").Value = 125018208
$ws.Range("L132").Value = 11156.25
$ws.Range("M132").Value = -125015678
$ws.Range("N132").Value = -16216.25

# ===================== Sheet LTW =====================
$ws = $wb.Worksheets.Item("LTW")
# row 16 (G16=5289)
$ws.Range("H16").Value = 820.8182
$ws.Range("I16").Value = 292.64706
$ws.Range("J16").Value = 2616.6
$ws.Range("K16").Value = 292.64706
$ws.Range("L16").Value = 2616.6
$ws.Range("M16").Value = -122.64706
$ws.Range("N16").Value = -2956.6
# row 68 (G68=12563)
$ws.Range("H68").Value = 2025
$ws.Range("I68").Value = 1366.6666
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 1366.6666
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -617.6666
$ws.Range("N68").Value = -5498
# row 71 (G71=12563)
$ws.Range("H71").Value = 2025
$ws.Range("I71").Value = 1366.6666
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 6833.333000000001
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -3089.333000000001
$ws.Range("N71").Value = -27488

# ===================== Sheet WVR =====================
$ws = $wb.Worksheets.Item("WVR")
# row 136 (G136=44031)
$ws.Range("H136").Value = 2697.077
$ws.Range("I136").Value = 2490.0645
$ws.Range("J136").Value = 3499.25
$ws.Range("K136").Value = 7470.193499999999
$ws.Range("L136").Value = 10497.75
$ws.Range("M136").Value = -4920.193499999999
$ws.Range("N136").Value = -15597.75
